# Fruta / hortaliza, semanal
# Insert a new weekly record at row 144 (pushing existing rows 144:170 down
# to 145:171) and populate it with this week's Uva price data for
# Terminal Hortofrutícola Agro Chillán.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 144:170 down to 145:171, leaving a blank row 144 to fill in.
$ws.Rows.Item(144).Insert()

$ws.Range("A144").Value = 7
$ws.Range("B144").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C144").Value = "Ñuble"
$ws.Range("D144").Value = 44995
$ws.Range("D144").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E144").Value = 16
$ws.Range("F144").Value = "Fruta"
$ws.Range("G144").Value = 100109
$ws.Range("H144").Value = "Uva"
$ws.Range("I144").Value = 100109001
$ws.Range("J144").Value = "Uva"
$ws.Range("K144").Value = "Flame Seedless"
$ws.Range("L144").Value = "Primera"
$ws.Range("M144").Value = 50
$ws.Range("N144").Value = 10000
$ws.Range("O144").Value = 10000
$ws.Range("P144").Value = 10000
$ws.Range("Q144").Value = "`$/bandeja 18 kilos"
$ws.Range("R144").Value = "Región de O'Higgins"
$ws.Range("S144").Value = 556
$ws.Range("T144").Value = 18
